$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(18, 1).Value = "'12/12/2025"
$ws.Cells.Item(18, 1).Style = "Normal"

$ws.Cells.Item(18, 2).Value = 12589.72
$ws.Cells.Item(18, 3).Value = 0.2026587542587263
$ws.Cells.Item(18, 4).Value = 0.7973412457412737
$ws.Cells.Item(18, 5).Value = -123.43
$ws.Cells.Item(18, 6).Value = -26.83
$ws.Cells.Item(18, 7).Value = -20478.46
$ws.Cells.Item(18, 8).Value = -67.11
$ws.Cells.Item(18, 9).Value = -351.43
$ws.Cells.Item(18, 10).Value = -12.11
